$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Must Haves" checklist (rows 9-16): all items now satisfied -> Yes (1)
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 1

# Bonus features (Protocol section): "Describes design" and "Contains link to GIT" now score points
$ws.Range("B47").Value = 1
$ws.Range("B52").Value = 0.5

# Refresh the "Sum Points" formula so its cached result reflects the new inputs
$ws.Range("B54").Formula = "=IF(MIN(B9:B16)=1,SUM(B21:B52),0)"

# Move the view / selection to where the author left off editing
$ws.Range("G45").Select() | Out-Null
